$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the species-record data between row 2 and row 3 for the
# columns that differ between the two observations (A, B, D, E, F, G, H, Q, R).
# Column C ("Valideringsstatus") is identical for both rows and is left alone.

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range("$col`2")
    $cell3 = $ws.Range("$col`3")

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
